$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (C, D, E, F, G) - A and B columns unchanged
$data = @{
    2  = @(48.05, 0.3075,  0.15375, 15.38142, 0.85)
    3  = @(31.35, 0.33127, 0.16564, 11.37465, 0.15)
    4  = @(20.2,  0.65003, 0.16251, 15.06344, 0.15)
    5  = @(26.25, 0.32031, 0.08008, 9.67203,  0.85)
    6  = @(19.1,  0.43267, 0.07211, 10.52543, 0.85)
    7  = @(13.65, 0.59325, 0.09887, 10.56027, 0.15)
    8  = @(9.800000000000001, 0.6511400000000001, 0.08139, 9.459720000000001, 0.15)
    9  = @(15.1,  0.40776, 0.05097, 8.48845,  0.85)
    10 = @(12,    0.42173, 0.04217, 7.52023,  0.85)
    11 = @(8.050000000000001, 0.55266, 0.05527, 7.25281,  0.15)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
